$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.397.91'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -4.86%  '
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.574.99'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  -4.63%  '
$ws.Range("E4").Value = '  +0.18%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.06%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '292.12'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -2.75%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3671'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -3.32%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.62'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -2.17%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3367'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -5.82%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.176'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -4.30%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07596'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -6.40%  '
$ws.Range("E12").Value = '  +0.17%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.23'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -4.14%  '
$ws.Range("E14").Value = '  -5.51%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.884'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -7.57%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001139'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -5.68%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.565.94'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -5.89%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.54'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -8.06%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06770'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -3.24%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.14%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.253'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -7.96%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.34'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -6.77%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.95'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -5.21%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.421.14'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -4.86%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.416'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -2.99%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.978'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +1.59%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.84'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -5.63%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '145.73'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -4.53%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.955'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  -5.31%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.24'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -6.15%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.745.73'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -5.48%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.292'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -10.02%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.980'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -7.42%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9744'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -5.80%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.50'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -11.55%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08522'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -2.44%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02545'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -7.31%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2297'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -6.50%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06549'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -5.03%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.510'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -8.33%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.84'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -11.11%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.263'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -4.63%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6397'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -7.90%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.64'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -6.86%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +0.07%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6016'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -6.94%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.781'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -3.81%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.130'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -6.50%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.78'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -5.12%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07284'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -6.92%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.184'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -0.90%  '
